# Weekly update: insert a new data row at row 46 (pushing the existing
# rows 46-149 down to 47-150) and populate the new row with the latest
# "Haba" price record for "Macroferia Regional de Talca".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 46:149 down to 47:150, creating a blank row 46.
$ws.Rows.Item(46).Insert()

# Populate the new row 46 with this week's record.
$ws.Range("A46").Value = 5
$ws.Range("B46").Value = "Macroferia Regional de Talca"
$ws.Range("C46").Value = "Maule"
$ws.Range("D46").Value = 45246
$ws.Range("E46").Value = 7
$ws.Range("F46").Value = 100112026
$ws.Range("G46").Value = "Haba"
$ws.Range("H46").Value = "Sin especificar"
$ws.Range("I46").Value = "Primera"
$ws.Range("J46").Value = 400
$ws.Range("K46").Value = 8000
$ws.Range("L46").Value = 8000
$ws.Range("M46").Value = 8000
$ws.Range("N46").Value = "`$/saco 25 kilos"
$ws.Range("O46").Value = "Región del Maule"
$ws.Range("P46").Value = 320
$ws.Range("Q46").Value = 25
$ws.Range("R46").Value = "Hortaliza"
